$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 82) to the daily log table.
# A82 must stay a literal text string ("2025/10/09"), not be auto-converted
# to a date serial number, so we temporarily force a text number format,
# assign the value, then clear formatting again so the cell ends up with
# the default (no explicit) style -- matching the rest of the data rows.
$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = "2025/10/09"
$ws.Range("B82").Value = "木"
$ws.Range("C82").Value = 7
$ws.Range("D82").Value = 136

$ws.Range("A82").ClearFormats()
